$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(5,5),
    @(7,7),
    @(9,9),
    @(2,3),
    @(7,7),
    @(3,5),
    @(8,9),
    @(4,5),
    @(6,7),
    @(6,6),
    @(7,7),
    @(8,9),
    @(8,8),
    @(7,8),
    @(8,9),
    @(10,10),
    @(6,7),
    @(5,6),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(1,2),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,4),
    @(4,5),
    @(3,4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
